# "adjusted comprehension questions for new version"
#
# The sheet previously had a duplicate pair of "attention" self-report
# questions (one about pictures, one about picture locations). The new
# version drops the second ("...a képek helyére...") variant, which
# shifts every following row up by one, and replaces the two old
# "how many pictures did you miss" follow-up questions with two new
# "which key did you press" comprehension questions. It also turns off
# wrapping on the long instruction cells in F2/F3 (their rows no longer
# need the huge fixed row height) and resets the sheet view/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old row 7 ("Mit gondol, ebben a körben milyen arányban
# tudott figyelni a képek helyére? ...") - everything below shifts up.
$ws.Rows("7:7").Delete()

# Old row 10 -> new row 9: replace the "first task - how many pictures
# did you not look at carefully" question with the new key-press question.
$ws.Range("A9").Value = "Melyik billentyűvel jelölte azokat a képeket, amelyeket beválogatott a galériába?"
$ws.Range("B9").Value = "D billentyű    F billentyű   J billentyű    K billentyű"
$ws.Range("C9").Value = "j"

# Old row 11 -> new row 10: replace the "second task - how many pictures
# did you not look at carefully" question with the new key-press question.
$ws.Range("A10").Value = "Melyik billentyűvel jelölte a képeket, amelyek pontosan ugyanolyanok voltak, mint a megelőző Galériaberendezés feladatban?"
$ws.Range("B10").Value = "D billentyű    F billentyű   J billentyű    K billentyű"
$ws.Range("C10").Value = "f"

# The long instruction cells F2/F3 no longer wrap (their rows shed the
# ht="409.5" fixed height in favor of the default row height).
$ws.Range("F2").WrapText = $false
$ws.Range("F3").WrapText = $false
$ws.Rows("2:2").AutoFit()
$ws.Rows("3:3").AutoFit()

# Reset the view: no more frozen/scrolled topLeftCell, selection moves to A6.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("A6").Select()
